$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45189 -> 45190, i.e. 2023-09-20 -> 2023-09-21) for every data row (rows 2-485).
$ws.Range("C2:C485").Value = 45190
